$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "tags" column (B) and the "tagsEn" column (original F).
$ws.Columns("F").Delete()
$ws.Columns("B").Delete()

# Update the selection to match the new layout.
$ws.Range("E6").Select()
